# Remove the now-obsolete "Syntax\" path segment from the CodeSnippets
# image paths stored in column A (rows 2-10) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value -like "*\Syntax\Snippets\*") {
        $cell.Value = $value -replace [regex]::Escape("\Syntax\Snippets\"), "\Snippets\"
    }
}
